# Re-curate the dimension/measure metadata columns for
# situacion-preferente, lugar-trabajo-o-estudio and tiempo-desplazamiento:
# they move from "iaest-dimension:*" / dim / skos:Concept (with an
# associated mapping-*.xlsx lookup file in row 5) to plain
# "iaest-measure:*" / medida / xsd:int columns - matching column A's
# personas-residentes-viviendas-familiares measure. The refArea column
# for "aragon" (E) is reclassified as a plain sdmx-dimension:refArea
# with a URI-Comunidad datatype instead of URI-Provincia/skos:Concept.
# The now-obsolete row 5 (mapping workbook references) is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - concept
$ws.Range("B2").Value = "iaest-measure:situacion-preferente"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "iaest-measure:lugar-trabajo-o-estudio"
$ws.Range("H2").Value = "iaest-measure:tiempo-desplazamiento"

# Row 3 - dim/medida
$ws.Range("B3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"

# Row 4 - datatype
$ws.Range("B4").Value = "xsd:int"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"

# Row 5 - drop the obsolete mapping-file references entirely
$ws.Rows("5:5").Delete()
